$d = $word.ActiveDocument
Write-Host "no-op"
